$d = $word.ActiveDocument

# --- 1. Swap the "Readings & Reference Material" / "Lecture: An Introduction
#        to Cartography" headings: the Lecture heading now comes first (and
#        is promoted from Heading 3 to Heading 2), followed by the Readings
#        & Reference Material heading (which stays Heading 2). ---

$p6 = $d.Paragraphs(6)
$p7 = $d.Paragraphs(7)

# Paragraph 6 currently: "Readings & Reference Material" (Heading 2)
# Paragraph 7 currently: "Lecture: An Introduction to Cartography" (Heading 3)

$p6.Range.Text = "Lecture: An Introduction to Cartography"
$p6.Style = "Heading 2"

$p7.Range.Text = "Readings & Reference Material"
$p7.Style = "Heading 2"

# --- 2. Remove the now-redundant "Reading" heading paragraph (paragraph 8)
#        entirely -- the reading list that used to sit under it now sits
#        directly under "Readings & Reference Material". ---

$p8 = $d.Paragraphs(8)
$p9 = $d.Paragraphs(9)
$killRange = $d.Range($p8.Range.Start, $p9.Range.Start)
$killRange.Delete()

Write-Output "done"
